$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the header row (row 1) with two more sequential values (14, 15)
# in columns P and Q, copying the existing header style (bold, bordered,
# centered) from O1 so the new cells match the rest of the header.
$ws.Range("O1").Copy($ws.Range("P1:Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Rows 2-25: the "1" markers that used to sit in columns I and M move two
# columns to the right (into K and O), and two new data columns P and Q
# are appended with value 2.
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
